$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = "'" + $value
    $ws.Cells.Item($row, $col).Style = "Normal"
}

# Row 2
Set-TextValue 2 4 "67.670.74"
Set-TextValue 2 5 "  +0.78%  "

# Row 3
Set-TextValue 3 4 "3.889.33"
Set-TextValue 3 5 "  +0.07%  "

# Row 4
Set-TextValue 4 4 "0.999"
Set-TextValue 4 5 "  -0.17%  "

# Row 5
Set-TextValue 5 4 "466.53"
Set-TextValue 5 5 "  +9.50%  "

# Row 6
Set-TextValue 6 4 "148.68"
Set-TextValue 6 5 "  +13.92%  "

# Row 7
Set-TextValue 7 4 "0.636"
Set-TextValue 7 5 "  +3.82%  "

# Row 8
Set-TextValue 8 4 "0.997"
Set-TextValue 8 5 "  -0.14%  "

# Row 9
Set-TextValue 9 4 "0.748"
Set-TextValue 9 5 "  +3.34%  "

# Row 10
Set-TextValue 10 5 "  -0.93%  "

# Row 11
Set-TextValue 11 5 "  -8.65%  "

# Row 12
Set-TextValue 12 4 "43.83"
Set-TextValue 12 5 "  +7.20%  "

# Row 13
Set-TextValue 13 4 "10.42"
Set-TextValue 13 5 "  +2.11%  "

# Row 14
Set-TextValue 14 4 "4.518.57"
Set-TextValue 14 5 "  +0.62%  "

# Row 15
Set-TextValue 15 4 "14.80"
Set-TextValue 15 5 "  -7.15%  "

# Row 16
Set-TextValue 16 4 "3.863.04"
Set-TextValue 16 5 "  -0.70%  "

# Row 17
Set-TextValue 17 5 "  -0.31%  "

# Row 18
Set-TextValue 18 4 "20.09"
Set-TextValue 18 5 "  +0.81%  "

# Row 19
Set-TextValue 19 5 "  +8.04%  "

# Row 20
Set-TextValue 20 4 "67.830.59"
Set-TextValue 20 5 "  +0.79%  "

# Row 21
Set-TextValue 21 4 "431.24"
Set-TextValue 21 5 "  +4.61%  "

# Row 22
Set-TextValue 22 4 "14.86"
Set-TextValue 22 5 "  +0.27%  "

# Row 23
Set-TextValue 23 5 "  +9.40%  "

# Row 24
Set-TextValue 24 4 "88.58"
Set-TextValue 24 5 "  +5.29%  "

# Row 25
Set-TextValue 25 4 "3.59"
Set-TextValue 25 5 "  +10.29%  "

# Row 26
$ws.Cells.Item(26, 2).Value = "EthereumClassic"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue 26 4 "37.92"
Set-TextValue 26 5 "  +1.16%  "

# Row 27
$ws.Cells.Item(27, 2).Value = "RenderToken"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue 27 4 "10.31"
Set-TextValue 27 5 "  +17.17%  "

# Row 28
Set-TextValue 28 4 "10.13"
Set-TextValue 28 5 "  +2.27%  "

# Row 29
Set-TextValue 29 4 "5.50"
Set-TextValue 29 5 "  +2.57%  "

# Row 30
Set-TextValue 30 4 "732.59"
Set-TextValue 30 5 "  -0.65%  "

# Row 31
Set-TextValue 31 5 "  +11.38%  "

# Row 32
Set-TextValue 32 4 "13.80"
Set-TextValue 32 5 "  +3.37%  "

# Row 33
Set-TextValue 33 5 "  -0.34%  "

# Row 34
Set-TextValue 34 4 "43.09"
Set-TextValue 34 5 "  +10.69%  "

# Row 35
Set-TextValue 35 4 "0.162"
Set-TextValue 35 5 "  +7.16%  "

# Row 36
Set-TextValue 36 4 "57.90"
Set-TextValue 36 5 "  +3.68%  "

# Row 37
Set-TextValue 37 5 "  +0.15%  "

# Row 38
Set-TextValue 38 4 "5.48"
Set-TextValue 38 5 "  +4.83%  "

# Row 39
Set-TextValue 39 4 "0.0480"
Set-TextValue 39 5 "  +4.26%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "ThetaToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue 40 4 "2.92"
Set-TextValue 40 5 "  +1.78%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "TheGraph"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue 41 4 "0.346"
Set-TextValue 41 5 "  +10.83%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "Stellar"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue 42 4 "0.142"
Set-TextValue 42 5 "  +5.87%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "PEPE"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue 43 4 "0.0₃0683"
Set-TextValue 43 5 "  -10.04%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "Fetch.AI"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue 44 4 "2.58"
Set-TextValue 44 5 "  +16.99%  "

# Row 45
Set-TextValue 45 5 "  -0.10%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "LidoDAOToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue 46 4 "3.44"
Set-TextValue 46 5 "  +2.47%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "WEMIXToken"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue 47 4 "2.78"
Set-TextValue 47 5 "  +8.19%  "

# Row 48
Set-TextValue 48 4 "3.24"
Set-TextValue 48 5 "  +1.34%  "

# Row 49
Set-TextValue 49 5 "  +5.17%  "

# Row 50
Set-TextValue 50 5 "  +3.44%  "

# Row 51
Set-TextValue 51 4 "144.48"
Set-TextValue 51 5 "  +1.89%  "
